# Changes of 4th May 2022
# CrudOperation.xlsx: rows 2-4 of Sheet1 had their "ExpectedRate" (column M)
# switched from a computed currency number to the literal text that was
# already used for the "ActualRate" column (Q), and the sheet's saved
# view state (scroll position / active selection) moved on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- M2:M4 — replace the numeric currency values with literal text ------
# Simply assigning a "$nn.nn"-shaped string to .Value lets Excel's smart
# entry re-parse it back into a currency number, so instead we go through
# a text formula (="$19.04") and then freeze it to a literal value with
# Copy/PasteSpecial(xlPasteValues). Resetting the cell style to "Normal"
# first drops the old currency number format/border styling the same way
# the saved workbook does.
$rateUpdates = @{
    "M2" = "$19.04"
    "M3" = "$17.98"
    "M4" = "$49.70"
}

foreach ($addr in $rateUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Style = "Normal"
    $cell.Formula = '="' + $rateUpdates[$addr] + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

$excel.CutCopyMode = $false

# --- sheet view state: scrolled one column left, selection moved --------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 10   # topLeftCell J1 (was K1)
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("P6").Select()                # was P4
